# Generated_TestCases.xlsx - "Add generated test cases from Claude"
#
# This script rewrites the "Testcases" worksheet:
#  - Replaces the component banner (B2) with a multi-line placeholder block
#    that also absorbs the former "MFP: Any" line (so E3 is cleared).
#  - Clears out the 16 detailed test-case rows (6-21), leaving each one
#    blank except for a "Not Executed" marker in column G (the Result
#    drop-down), matching a freshly generated / not-yet-run template.
#  - Clears the performance-test row (22) and the usability-test row (23)
#    entirely (row 23 ends up with no content at all).
#  - Converts the old "TC019" summary row (24) into a "Test Summary" label
#    and the old "TC020" row (25) into a "Test Case Count:" label.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Testcases")

# --- Header block -----------------------------------------------------
$ws.Range("B2").Value = "Component: MultiFunctionalTool`nMFP: Any`nBuild: `nDate: `nTarget: "
$ws.Range("E3").ClearContents()

# --- Detailed test case rows (6-21): wipe details, mark Not Executed --
for ($r = 6; $r -le 21; $r++) {
    $ws.Range("B${r}:F${r}").ClearContents()
    $ws.Range("G${r}").Value = "Not Executed"
    $ws.Range("H${r}").ClearContents()
}

# --- Row 22 (performance/response-time case): wipe entirely -----------
$ws.Range("A22:H22").ClearContents()

# --- Row 23 (usability case): wipe entirely (no styled cells remain) --
$ws.Range("A23:H23").ClearContents()

# --- Row 24: becomes the "Test Summary" section header -----------------
$ws.Range("B24").Value = "Test Summary"
$ws.Range("C24").ClearContents()
$ws.Range("D24:H24").ClearContents()

# --- Row 25: becomes the "Test Case Count:" label ----------------------
$ws.Range("B25:C25").ClearContents()
$ws.Range("D25").Value = "Test Case Count:"
$ws.Range("E25:H25").ClearContents()
